# Adding Mac and Windows Icons
# Duplicate the existing "Mac" hexagon/E slide (slide 1) to create the new
# "Windows" variant, then nudge its geometry and swap the run font so it
# matches the Bevan/Windows icon. PowerPoint inserts the duplicate right
# after the source slide, which is exactly where the new slide belongs.

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$sourceSlide = $p.Slides.Item(1)

$newRange = $sourceSlide.Duplicate()
$newSlide = $newRange.Item(1)

# Shape 1: the dark hexagon background - only its vertical position moves.
$hexagon = $newSlide.Shapes.Item(1)
$hexagon.Top = (93306 + 0.5) / $EMU_PER_PT

# Shape 2: the big "E" textbox - vertical position + height change, and the
# run switches from the Mac "Berlin Sans FB Demi" face to the Windows
# "Bevan" face at a slightly smaller point size.
$textBox = $newSlide.Shapes.Item(2)
$textBox.Top = (788082 + 0.5) / $EMU_PER_PT
$textBox.Height = (7017306 + 0.5) / $EMU_PER_PT

$runRange = $textBox.TextFrame.TextRange
$runRange.Font.Size = 450
$runRange.Font.Name = "Bevan"
